$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy row 6 cell formatting into the new rows (7-13)
$ws.Range("A6:BF6").Copy() | Out-Null
$ws.Range("A7:BF13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Row 7 ----
$ws.Range("A7").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B7").Value = "segmentectomia 5"
$ws.Range("C7").Value = "Segmentectomia1a8"
$ws.Range("D7").Value = 1676
$ws.Range("E7").Value = "19/07/2019"
$ws.Range("F7").Value = "19/07/2019"
$ws.Range("G7").Value = "26/6/19"
$ws.Range("H7").Value = "V,IV"
$ws.Range("I7").Value = 1534
$ws.Range("J7").Value = "Francesc"
$ws.Range("K7").Value = "Dolera"
$ws.Range("L7").Value = "bernal"
$ws.Range("M7").Value = 12817360
$ws.Range("N7").Value = 43725
$ws.Range("O7").Value = "Si"
$ws.Range("P7").Value = "Home"
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "60"
$ws.Range("R7").NumberFormat = "@"
$ws.Range("R7").Value = "80"
$ws.Range("S7").Value = 170
$ws.Range("T7").Value = 26
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = "No"
$ws.Range("W7").Value = "Resecció Menor (<3 segm)"
$ws.Range("X7").Value = "Oberta"
$ws.Range("Y7").Value = "No"
$ws.Range("Z7").Value = "No"
$ws.Range("AA7").Value = "Impressió R0"
$ws.Range("AB7").Value = 1
$ws.Range("AC7").Value = 2
$ws.Range("AD7").Value = "No"
$ws.Range("AE7").Value = "No"
$ws.Range("AF7").Value = "No"
$ws.Range("AG7").Value = "No"
$ws.Range("AH7").NumberFormat = "@"
$ws.Range("AH7").Value = "0"
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = 1
$ws.Range("AK7").Value = 5
$ws.Range("AL7").Value = 1
$ws.Range("AM7").Value = "No"
$ws.Range("AO7").Value = 44421
$ws.Range("AP7").Value = "No"
$ws.Range("AQ7").Value = "No"
$ws.Range("AR7").Value = "Viu"
$ws.Range("AS7").Value = "No"
$ws.Range("AT7").Value = "No"
$ws.Range("AU7").Value = "No"
$ws.Range("AV7").Value = "No"
$ws.Range("AW7").Value = "No"
$ws.Range("AX7").Value = "No"
$ws.Range("AY7").Value = "No"
$ws.Range("BB7").Value = 6
$ws.Range("BC7").Value = 44614.35327677083
$ws.Range("BD7").Value = 43175
$ws.Range("BE7").Value = "No"

# ---- Row 8 ----
$ws.Range("A8").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B8").Value = "segmentectomia 5"
$ws.Range("C8").Value = "Segmentectomia1a8"
$ws.Range("D8").Value = 1676
$ws.Range("E8").Value = "19/07/2019"
$ws.Range("F8").Value = "19/07/2019"
$ws.Range("G8").Value = "26/6/19"
$ws.Range("H8").Value = "V,IV"
$ws.Range("I8").Value = 1534
$ws.Range("J8").Value = "Francesc"
$ws.Range("K8").Value = "Dolera"
$ws.Range("L8").Value = "bernal"
$ws.Range("M8").Value = 12817360
$ws.Range("N8").Value = 43725
$ws.Range("O8").Value = "Si"
$ws.Range("P8").Value = "Home"
$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = "60"
$ws.Range("R8").NumberFormat = "@"
$ws.Range("R8").Value = "80"
$ws.Range("S8").Value = 170
$ws.Range("T8").Value = 26
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = "No"
$ws.Range("W8").Value = "Resecció Menor (<3 segm)"
$ws.Range("X8").Value = "Oberta"
$ws.Range("Y8").Value = "No"
$ws.Range("Z8").Value = "No"
$ws.Range("AA8").Value = "Impressió R0"
$ws.Range("AB8").Value = 1
$ws.Range("AC8").Value = 2
$ws.Range("AD8").Value = "No"
$ws.Range("AE8").Value = "No"
$ws.Range("AF8").Value = "No"
$ws.Range("AG8").Value = "No"
$ws.Range("AH8").NumberFormat = "@"
$ws.Range("AH8").Value = "0"
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = 1
$ws.Range("AK8").Value = 5
$ws.Range("AL8").Value = 1
$ws.Range("AM8").Value = "No"
$ws.Range("AO8").Value = 44421
$ws.Range("AP8").Value = "No"
$ws.Range("AQ8").Value = "No"
$ws.Range("AR8").Value = "Viu"
$ws.Range("AS8").Value = "No"
$ws.Range("AT8").Value = "No"
$ws.Range("AU8").Value = "No"
$ws.Range("AV8").Value = "No"
$ws.Range("AW8").Value = "No"
$ws.Range("AX8").Value = "No"
$ws.Range("AY8").Value = "No"
$ws.Range("BB8").Value = 6
$ws.Range("BC8").Value = 44614.35377028935
$ws.Range("BD8").Value = 43175
$ws.Range("BE8").Value = "No"

# ---- Row 9 ----
$ws.Range("A9").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B9").Value = "segmentect 5 i 2 RL"
$ws.Range("C9").Value = "Segmentectomia1a8"
$ws.Range("D9").Value = 1694
$ws.Range("E9").Value = "18/10/2019"
$ws.Range("F9").Value = "18/10/2019"
$ws.Range("G9").Value = "16/07/2019"
$ws.Range("H9").Value = "II,V,VIII,VII"
$ws.Range("I9").Value = 1548
$ws.Range("J9").Value = "Fco jesus"
$ws.Range("K9").Value = "Herrera"
$ws.Range("L9").Value = "Exposito"
$ws.Range("M9").Value = 13611095
$ws.Range("N9").Value = 43846
$ws.Range("O9").Value = "Si"
$ws.Range("P9").Value = "Home"
$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "61"
$ws.Range("R9").NumberFormat = "@"
$ws.Range("R9").Value = "76"
$ws.Range("S9").Value = 175
$ws.Range("T9").Value = 25
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = "No"
$ws.Range("W9").Value = "Resecció Menor (<3 segm)"
$ws.Range("X9").Value = "1er temps (mobilització)"
$ws.Range("Y9").Value = "No"
$ws.Range("Z9").Value = "No"
$ws.Range("AA9").Value = "Impressió R1"
$ws.Range("AB9").Value = 3
$ws.Range("AC9").Value = 2
$ws.Range("AD9").Value = "Si"
$ws.Range("AE9").Value = "No"
$ws.Range("AF9").Value = "No"
$ws.Range("AG9").Value = "Si"
$ws.Range("AH9").Value = "IIIa"
$ws.Range("AJ9").Value = 3
$ws.Range("AK9").Value = 2.5
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = "Si"
$ws.Range("AN9").Value = "ampli quirúgica i Aquamantis"
$ws.Range("AO9").Value = 44243
$ws.Range("AP9").Value = "No"
$ws.Range("AQ9").Value = "No"
$ws.Range("AR9").Value = "Viu"
$ws.Range("AS9").Value = "No"
$ws.Range("AT9").Value = "Si"
$ws.Range("AU9").Value = "No"
$ws.Range("AV9").Value = "No"
$ws.Range("AW9").Value = "No"
$ws.Range("AX9").Value = "No"
$ws.Range("AY9").Value = "Si"
$ws.Range("AZ9").Value = "Si"
$ws.Range("BA9").Value = "Si"
$ws.Range("BB9").Value = 6
$ws.Range("BC9").Value = 44614.35471002315
$ws.Range("BD9").Value = 43721

# ---- Row 10 ----
$ws.Range("A10").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B10").Value = "segmentect 5 i 2 RL"
$ws.Range("C10").Value = "Segmentectomia1a8"
$ws.Range("D10").Value = 1694
$ws.Range("E10").Value = "18/10/2019"
$ws.Range("F10").Value = "18/10/2019"
$ws.Range("G10").Value = "16/07/2019"
$ws.Range("H10").Value = "II,V,VIII,VII"
$ws.Range("I10").Value = 1548
$ws.Range("J10").Value = "Fco jesus"
$ws.Range("K10").Value = "Herrera"
$ws.Range("L10").Value = "Exposito"
$ws.Range("M10").Value = 13611095
$ws.Range("N10").Value = 43846
$ws.Range("O10").Value = "Si"
$ws.Range("P10").Value = "Home"
$ws.Range("Q10").NumberFormat = "@"
$ws.Range("Q10").Value = "61"
$ws.Range("R10").NumberFormat = "@"
$ws.Range("R10").Value = "76"
$ws.Range("S10").Value = 175
$ws.Range("T10").Value = 25
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = "No"
$ws.Range("W10").Value = "Resecció Menor (<3 segm)"
$ws.Range("X10").Value = "1er temps (mobilització)"
$ws.Range("Y10").Value = "No"
$ws.Range("Z10").Value = "No"
$ws.Range("AA10").Value = "Impressió R1"
$ws.Range("AB10").Value = 3
$ws.Range("AC10").Value = 2
$ws.Range("AD10").Value = "Si"
$ws.Range("AE10").Value = "No"
$ws.Range("AF10").Value = "No"
$ws.Range("AG10").Value = "Si"
$ws.Range("AH10").Value = "IIIa"
$ws.Range("AJ10").Value = 3
$ws.Range("AK10").Value = 2.5
$ws.Range("AL10").Value = 0
$ws.Range("AM10").Value = "Si"
$ws.Range("AN10").Value = "ampli quirúgica i Aquamantis"
$ws.Range("AO10").Value = 44243
$ws.Range("AP10").Value = "No"
$ws.Range("AQ10").Value = "No"
$ws.Range("AR10").Value = "Viu"
$ws.Range("AS10").Value = "No"
$ws.Range("AT10").Value = "Si"
$ws.Range("AU10").Value = "No"
$ws.Range("AV10").Value = "No"
$ws.Range("AW10").Value = "No"
$ws.Range("AX10").Value = "No"
$ws.Range("AY10").Value = "Si"
$ws.Range("AZ10").Value = "Si"
$ws.Range("BA10").Value = "Si"
$ws.Range("BB10").Value = 6
$ws.Range("BC10").Value = 44614.35498545139
$ws.Range("BD10").Value = 43721

# ---- Row 11 ----
$ws.Range("A11").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B11").Value = "segmentect 5 i 2 RL"
$ws.Range("C11").Value = "Segmentectomia1a8"
$ws.Range("D11").Value = 1694
$ws.Range("E11").Value = "18/10/2019"
$ws.Range("F11").Value = "18/10/2019"
$ws.Range("G11").Value = "16/07/2019"
$ws.Range("H11").Value = "II,V,VIII,VII"
$ws.Range("I11").Value = 1548
$ws.Range("J11").Value = "Fco jesus"
$ws.Range("K11").Value = "Herrera"
$ws.Range("L11").Value = "Exposito"
$ws.Range("M11").Value = 13611095
$ws.Range("N11").Value = 43846
$ws.Range("O11").Value = "Si"
$ws.Range("P11").Value = "Home"
$ws.Range("Q11").NumberFormat = "@"
$ws.Range("Q11").Value = "61"
$ws.Range("R11").NumberFormat = "@"
$ws.Range("R11").Value = "76"
$ws.Range("S11").Value = 175
$ws.Range("T11").Value = 25
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = "No"
$ws.Range("W11").Value = "Resecció Menor (<3 segm)"
$ws.Range("X11").Value = "1er temps (mobilització)"
$ws.Range("Y11").Value = "No"
$ws.Range("Z11").Value = "No"
$ws.Range("AA11").Value = "Impressió R1"
$ws.Range("AB11").Value = 3
$ws.Range("AC11").Value = 2
$ws.Range("AD11").Value = "Si"
$ws.Range("AE11").Value = "No"
$ws.Range("AF11").Value = "No"
$ws.Range("AG11").Value = "Si"
$ws.Range("AH11").Value = "IIIa"
$ws.Range("AJ11").Value = 3
$ws.Range("AK11").Value = 2.5
$ws.Range("AL11").Value = 0
$ws.Range("AM11").Value = "Si"
$ws.Range("AN11").Value = "ampli quirúgica i Aquamantis"
$ws.Range("AO11").Value = 44243
$ws.Range("AP11").Value = "No"
$ws.Range("AQ11").Value = "No"
$ws.Range("AR11").Value = "Viu"
$ws.Range("AS11").Value = "No"
$ws.Range("AT11").Value = "Si"
$ws.Range("AU11").Value = "No"
$ws.Range("AV11").Value = "No"
$ws.Range("AW11").Value = "No"
$ws.Range("AX11").Value = "No"
$ws.Range("AY11").Value = "Si"
$ws.Range("AZ11").Value = "Si"
$ws.Range("BA11").Value = "Si"
$ws.Range("BB11").Value = 6
$ws.Range("BC11").Value = 44614.35590646991
$ws.Range("BD11").Value = 43721

# ---- Row 12 ----
$ws.Range("A12").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B12").Value = "segmentect 5 i 2 RL"
$ws.Range("C12").Value = "Segmentectomia1a8"
$ws.Range("D12").Value = 1694
$ws.Range("E12").Value = "18/10/2019"
$ws.Range("F12").Value = "18/10/2019"
$ws.Range("G12").Value = "16/07/2019"
$ws.Range("H12").Value = "II,V,VIII,VII"
$ws.Range("I12").Value = 1548
$ws.Range("J12").Value = "Fco jesus"
$ws.Range("K12").Value = "Herrera"
$ws.Range("L12").Value = "Exposito"
$ws.Range("M12").Value = 13611095
$ws.Range("N12").Value = 43846
$ws.Range("O12").Value = "Si"
$ws.Range("P12").Value = "Home"
$ws.Range("Q12").NumberFormat = "@"
$ws.Range("Q12").Value = "61"
$ws.Range("R12").NumberFormat = "@"
$ws.Range("R12").Value = "76"
$ws.Range("S12").Value = 175
$ws.Range("T12").Value = 25
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = "No"
$ws.Range("W12").Value = "Resecció Menor (<3 segm)"
$ws.Range("X12").Value = "1er temps (mobilització)"
$ws.Range("Y12").Value = "No"
$ws.Range("Z12").Value = "No"
$ws.Range("AA12").Value = "Impressió R1"
$ws.Range("AB12").Value = 3
$ws.Range("AC12").Value = 2
$ws.Range("AD12").Value = "Si"
$ws.Range("AE12").Value = "No"
$ws.Range("AF12").Value = "No"
$ws.Range("AG12").Value = "Si"
$ws.Range("AH12").Value = "IIIa"
$ws.Range("AJ12").Value = 3
$ws.Range("AK12").Value = 2.5
$ws.Range("AL12").Value = 0
$ws.Range("AM12").Value = "Si"
$ws.Range("AN12").Value = "ampli quirúgica i Aquamantis"
$ws.Range("AO12").Value = 44243
$ws.Range("AP12").Value = "No"
$ws.Range("AQ12").Value = "No"
$ws.Range("AR12").Value = "Viu"
$ws.Range("AS12").Value = "No"
$ws.Range("AT12").Value = "Si"
$ws.Range("AU12").Value = "No"
$ws.Range("AV12").Value = "No"
$ws.Range("AW12").Value = "No"
$ws.Range("AX12").Value = "No"
$ws.Range("AY12").Value = "Si"
$ws.Range("AZ12").Value = "Si"
$ws.Range("BA12").Value = "Si"
$ws.Range("BB12").Value = 6
$ws.Range("BC12").Value = 44614.35669915509
$ws.Range("BD12").Value = 43721

# ---- Row 13 ----
$ws.Range("A13").Value = "Segmentectomia o Bisegmentectomia"
$ws.Range("B13").Value = "segmentect 5 i 2 RL"
$ws.Range("C13").Value = "Segmentectomia1a8"
$ws.Range("D13").Value = 1694
$ws.Range("E13").Value = "18/10/2019"
$ws.Range("F13").Value = "18/10/2019"
$ws.Range("G13").Value = "16/07/2019"
$ws.Range("H13").Value = "II,V,VIII,VII"
$ws.Range("I13").Value = 1548
$ws.Range("J13").Value = "Fco jesus"
$ws.Range("K13").Value = "Herrera"
$ws.Range("L13").Value = "Exposito"
$ws.Range("M13").Value = 13611095
$ws.Range("N13").Value = 43846
$ws.Range("O13").Value = "Si"
$ws.Range("P13").Value = "Home"
$ws.Range("Q13").NumberFormat = "@"
$ws.Range("Q13").Value = "61"
$ws.Range("R13").NumberFormat = "@"
$ws.Range("R13").Value = "76"
$ws.Range("S13").Value = 175
$ws.Range("T13").Value = 25
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = "No"
$ws.Range("W13").Value = "Resecció Menor (<3 segm)"
$ws.Range("X13").Value = "1er temps (mobilització)"
$ws.Range("Y13").Value = "No"
$ws.Range("Z13").Value = "No"
$ws.Range("AA13").Value = "Impressió R1"
$ws.Range("AB13").Value = 3
$ws.Range("AC13").Value = 2
$ws.Range("AD13").Value = "Si"
$ws.Range("AE13").Value = "No"
$ws.Range("AF13").Value = "No"
$ws.Range("AG13").Value = "Si"
$ws.Range("AH13").Value = "IIIa"
$ws.Range("AJ13").Value = 3
$ws.Range("AK13").Value = 2.5
$ws.Range("AL13").Value = 0
$ws.Range("AM13").Value = "Si"
$ws.Range("AN13").Value = "ampli quirúgica i Aquamantis"
$ws.Range("AO13").Value = 44243
$ws.Range("AP13").Value = "No"
$ws.Range("AQ13").Value = "No"
$ws.Range("AR13").Value = "Viu"
$ws.Range("AS13").Value = "No"
$ws.Range("AT13").Value = "Si"
$ws.Range("AU13").Value = "No"
$ws.Range("AV13").Value = "No"
$ws.Range("AW13").Value = "No"
$ws.Range("AX13").Value = "No"
$ws.Range("AY13").Value = "Si"
$ws.Range("AZ13").Value = "Si"
$ws.Range("BA13").Value = "Si"
$ws.Range("BB13").Value = 6
$ws.Range("BC13").Value = 44614.35753907407
$ws.Range("BD13").Value = 43721

# Step 3: restore exact per-cell formatting (in case any Value assignment drifted number formats)
$ws.Range("A6:BF6").Copy() | Out-Null
$ws.Range("A7:BF7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:BF8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:BF9").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:BF10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:BF11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:BF12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:BF13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
